$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 13 ("Programa resumido:" row's
# value row), pushing the current rows 13-23 down to 14-24.
$ws.Rows.Item(13).Insert()

# The insert auto-extends column A's styling into the new row 13; the
# target layout has no A13 cell at all, so drop it completely.
$ws.Range("A13").Clear()

# New row 13 holds the "Docentes responsáveis:" value (previously
# mis-placed at A10/B10/C10 as part of the "Objetivos:" row).
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "144651 - Antonio Fernando Sartori"
$ws.Range("C13").Value = "144651 - Antonio Fernando Sartori"

# Row 10 ("Objetivos:") gets its real Portuguese objectives text instead
# of the misplaced professor name.
$ws.Range("B10").Value = "A disciplina visa propiciar aos alunos os conhecimentos básicos de eletroquímica, tanto do ponto de vista da eletroquímica iônica como da eletródica, e apresentar as principais aplicações da eletroquímica"
$ws.Range("C10").Value = "A disciplina visa propiciar aos alunos os conhecimentos básicos de eletroquímica, tanto do ponto de vista da eletroquímica iônica como da eletródica, e apresentar as principais aplicações da eletroquímica"

# Row 14 ("Programa resumido:") gets the real short-syllabus text instead
# of the placeholder "Semestral".
$ws.Range("B14").Value = "Princípios da eletroquímica iônica e da eletroquímica eletródica. Aplicações."
$ws.Range("C14").Value = "Princípios da eletroquímica iônica e da eletroquímica eletródica. Aplicações."

# Row 16 ("Programa:") gets the real full-syllabus text instead of the
# misplaced activation date.
$ws.Range("B16").Value = "Princípios da eletroquímica iônica: interações iônicas, equilíbrio iônico e condução eletrolítica. Princípios da eletroquímica eletródica: fenômenos interfaciais, potenciais de eletrodo e células eletroquímicas. Processos de eletrodo. Métodos eletroquímicos de análise química. Aplicações da eletroquímica: fontes eletroquímicas de energia, processos eletrometalúrgicos e galvanoplastia."
$ws.Range("C16").Value = "Princípios da eletroquímica iônica: interações iônicas, equilíbrio iônico e condução eletrolítica. Princípios da eletroquímica eletródica: fenômenos interfaciais, potenciais de eletrodo e células eletroquímicas. Processos de eletrodo. Métodos eletroquímicos de análise química. Aplicações da eletroquímica: fontes eletroquímicas de energia, processos eletrometalúrgicos e galvanoplastia."

# Row 19 ("Método:") gets the real teaching-method text instead of the
# misplaced professor name.
$ws.Range("B19").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Range("C19").Value = "Aulas expositivas, seminários e exercícios comentados."

# Row 20 ("Critério:") gets the real grading-criteria text instead of the
# teaching-method text.
$ws.Range("B20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Range("C20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."

# Row 21 ("Norma de recuperação:") gets the real make-up-exam rule text
# instead of the grading-criteria text.
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

# Row 22 ("Bibliografia:") gets the real bibliography text instead of the
# make-up-exam rule text.
$ws.Range("B22").Value = "BOCKRIS, J.O.M.;. REDDY, A.K.N. Modern Electrochemistry, 2 vols, Plenum Press, NY, 1977. DENARO, A. R. Fundamentos de Eletroquímica, Ed. Edgard Blucher, São Paulo, 1974. OLDHAM, K. B.; MYLAND, J. C. Fundamentals of Electrochemical Science, Academic Press, New York, 1994. TICIANELLI, E. A.; GONZALEZ, E. R., Eletroquímica, Edusp, 1998. CROW, D.R. Principles and Applications of Electrochemistry, Blackie Academic and Professional, London, 1994. KUHN, A .T. Industrial Electrochemical Processes, Elsevier, Amsterdam, 1971. PLETCHER, D.; WALSH, F. C. Industrial Electrochemistry, 2 ed., Blackie Academic & Professional, Cambridge,1993."
$ws.Range("C22").Value = "BOCKRIS, J.O.M.;. REDDY, A.K.N. Modern Electrochemistry, 2 vols, Plenum Press, NY, 1977. DENARO, A. R. Fundamentos de Eletroquímica, Ed. Edgard Blucher, São Paulo, 1974. OLDHAM, K. B.; MYLAND, J. C. Fundamentals of Electrochemical Science, Academic Press, New York, 1994. TICIANELLI, E. A.; GONZALEZ, E. R., Eletroquímica, Edusp, 1998. CROW, D.R. Principles and Applications of Electrochemistry, Blackie Academic and Professional, London, 1994. KUHN, A .T. Industrial Electrochemical Processes, Elsevier, Amsterdam, 1971. PLETCHER, D.; WALSH, F. C. Industrial Electrochemistry, 2 ed., Blackie Academic & Professional, Cambridge,1993."
